$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-13 from
# 45204 (2023-10-05) to 45205 (2023-10-06).
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
